$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (bold, border, center/top alignment) from the existing
# "sum" header (G1) onto the new "Save" header cell (H1), then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Populate the new "Save" column values for rows 2-7
$saveValues = @(1, 0, 1, 0, 0, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
